$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove "(sin resolver)" suffix from several comment strings ---
$ws.Range("G4").Value = "Recuerda que si el cliente despues de tres intentos no contesta agregar esa información"
$ws.Range("G5").Value = "En plan de riesgos falta efecto de la causa no conexión por admiadmin impide conectar a maquina cliente por ejemplo."
$ws.Range("G6").Value = "Al no requerir implementacion no debe decir los mensajes de implementacion"

# --- Row 7: status changes from "En proceso" to "Cerrada"; comment loses "(sin resolver)" suffix ---
$ws.Range("F7").Value = "Cerrada"
$ws.Range("G7").Value = "En la parte equipo de empresa poner SOS Software y en cliente el nombre de cliente"

# --- Row 8: clear out the stray merged-in row, keep only A8/C8 ---
$ws.Range("B8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("G8").ClearContents()

# --- Row heights ---
$ws.Rows.Item(5).RowHeight = 55.2
$ws.Rows.Item(7).RowHeight = 41.75
$ws.Rows.Item(8).RowHeight = 13.8

# --- Sheet view: reset scroll/selection back to A1 ---
[void]$ws.Range("A1").Select()

Write-Host "done"
